$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two rows near the top for "Preface" and "Summary" (new rows 2-3) ---
$ws.Rows("2:3").Insert()

# --- 2. Insert two rows before "Total" for "Bibliography" and "Appendix" ---
#     (at this point, "Total" lives at row 12, so inserting at 12:13 pushes it to row 14)
$ws.Rows("12:13").Insert()

# --- 3. Column A labels for the new rows ---
$ws.Range("A2").Value = "Preface"
$ws.Range("A3").Value = "Summary"
$ws.Range("A12").Value = "Bibliography"
$ws.Range("A13").Value = "Appendix"
$ws.Range("A15").Value = "Extra's"

# --- 4. Column B (Pages) for the new rows ---
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0

# --- 5. Fix up existing B-column values that changed (Project Plan 1.5 -> 1) ---
$ws.Range("B5").Value = 1

# --- 6. Headers for the two new columns ---
$ws.Range("C1").Value = "Completion"
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Value = "Tasks"
$ws.Range("D1").Font.Bold = $true

# --- 7. Column C (Completion %) values, rows 2-13 ---
$ws.Range("C2:C13").NumberFormat = "0%"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0.8
$ws.Range("C8").Value = 0.8
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0.5
$ws.Range("C13").Value = 0

# --- 8. Column D (Tasks) notes ---
$ws.Range("D7").Value = "review coupling part"
$ws.Range("D8").Value = "review?"
$ws.Range("D10").Value = "add lqr results (check for lqr + initial)`nadd result comparison"
$ws.Range("D10").WrapText = $true
$ws.Range("D12").Value = "fix missing links"
$ws.Range("D13").Value = "add models"
$ws.Range("D15").Value = "check labels`ncheck consistency"
$ws.Range("D15").WrapText = $true

# --- 9. Total row (now row 14): fix the SUM/AVERAGE ranges for the expanded table ---
$ws.Range("B14").Formula = "=SUM(B2:B13)"
$ws.Range("C14").Formula = "=AVERAGE(C2:C13)"
$ws.Range("C14").NumberFormat = "0%"

# --- 10. Row heights for wrapped rows ---
$ws.Rows("10").RowHeight = 28.8
$ws.Rows("15").RowHeight = 28.8

# --- 11. Column widths for the new columns ---
$ws.Columns("C").ColumnWidth = 12.88671875
$ws.Columns("D").ColumnWidth = 35.6640625

# --- 12. Selection matches the saved workbook state ---
$ws.Range("C17").Select()
